$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("données")

# Delete rows 17 and 18 (the rows with # pièce "19816" and "92317"),
# which shifts all subsequent rows up by two.
$ws.Range("A17:L18").EntireRow.Delete()

# Apply an AutoFilter over the data range (now A1:L78 after the row deletion).
$ws.Range("A1:L78").AutoFilter() | Out-Null

# AutoFilter registers a hidden sheet-scoped _FilterDatabase defined name.
$filterName = $ws.Names.Add("_xlnm._FilterDatabase", "=données!`$A`$1:`$L`$78")
$filterName.Visible = $false

# Move the selection, matching where the user's cursor ended up after editing.
$ws.Range("C19").Select() | Out-Null
